$d = $word.ActiveDocument

# Update the date line at the top of the document.
$d.Content.Find.Execute("2025-01-25 Saturday", $true, $false, $false, $false, $false,
                         $true, 1, $false, "2025-01-26 Sunday", 2)

# Update the 25 multiplication equations in the table, addressed by
# row/column so that identical text values occurring more than once
# (e.g. "559x7=3913" appears both as an old and a new value) cannot
# collide with each other the way a global find/replace would.
$t = $d.Tables.Item(1)

$t.Cell(1,1).Range.Text  = "265×6=1590"
$t.Cell(1,2).Range.Text  = "426×7=2982"
$t.Cell(1,3).Range.Text  = "915×5=4575"
$t.Cell(1,4).Range.Text  = "107×2=214"
$t.Cell(1,5).Range.Text  = "273×2=546"

$t.Cell(5,1).Range.Text  = "173×9=1557"
$t.Cell(5,2).Range.Text  = "605×5=3025"
$t.Cell(5,3).Range.Text  = "274×4=1096"
$t.Cell(5,4).Range.Text  = "708×2=1416"
$t.Cell(5,5).Range.Text  = "426×2=852"

$t.Cell(10,1).Range.Text = "232×6=1392"
$t.Cell(10,2).Range.Text = "526×4=2104"
$t.Cell(10,3).Range.Text = "548×5=2740"
$t.Cell(10,4).Range.Text = "778×2=1556"
$t.Cell(10,5).Range.Text = "559×7=3913"

$t.Cell(15,1).Range.Text = "803×7=5621"
$t.Cell(15,2).Range.Text = "851×4=3404"
$t.Cell(15,3).Range.Text = "123×6=738"
$t.Cell(15,4).Range.Text = "457×7=3199"
$t.Cell(15,5).Range.Text = "419×3=1257"

$t.Cell(20,1).Range.Text = "809×9=7281"
$t.Cell(20,2).Range.Text = "887×4=3548"
$t.Cell(20,3).Range.Text = "895×4=3580"
$t.Cell(20,4).Range.Text = "961×9=8649"
$t.Cell(20,5).Range.Text = "675×6=4050"
